# feat(firestarter): added argument --config. Played with the colors
#
# The underlying data cleanup: clear out a batch of stale MSCI_World (col B)
# and MSCI_ACWI_IMI (col D) values that shouldn't have been populated for
# these months, across the "Performance Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (MSCI_World) cleanups
$ws.Range("B72:B87").ClearContents()
$ws.Range("B136:B143").ClearContents()
$ws.Range("B268:B273").ClearContents()
$ws.Range("B301:B305").ClearContents()

# Column D (MSCI_ACWI_IMI) cleanups
$ws.Range("D25:D35").ClearContents()
$ws.Range("D138:D142").ClearContents()
$ws.Range("D188:D196").ClearContents()
$ws.Range("D274:D279").ClearContents()
$ws.Range("D305:D321").ClearContents()

# Update the view state to match where the user ended up scrolling/selecting.
$excel.Goto($ws.Range("A298"), $true)
$ws.Range("F315").Select()
